$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5675.4165
$ws.Range("I43").Value = 6332.5557
$ws.Range("J43").Value = 3704
$ws.Range("K43").Value = 6332.5557
$ws.Range("L43").Value = 3704
$ws.Range("M43").Value = -6263.5557
$ws.Range("N43").Value = -3842

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7443.1333
$ws.Range("I51").Value = 14824.5
$ws.Range("J51").Value = 6307.5386
$ws.Range("K51").Value = 14824.5
$ws.Range("L51").Value = 6307.5386
$ws.Range("M51").Value = -14340.5
$ws.Range("N51").Value = -7275.5386

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3848.9443
$ws.Range("I74").Value = 3163.5881
$ws.Range("K74").Value = 3163.5881
$ws.Range("M74").Value = -2227.5881

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3848.9443
$ws.Range("I77").Value = 3163.5881
$ws.Range("K77").Value = 15817.9405
$ws.Range("M77").Value = -11137.9405

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2069.152
$ws.Range("I132").Value = 1739.5128
$ws.Range("J132").Value = 3905.7144
$ws.Range("K132").Value = 5218.538399999999
$ws.Range("L132").Value = 11717.1432
$ws.Range("M132").Value = -2688.538399999999
$ws.Range("N132").Value = -16777.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 916.63635
$ws.Range("I2").Value = 958.5833
$ws.Range("J2").Value = 866.3
$ws.Range("K2").Value = 958.5833
$ws.Range("L2").Value = 866.3
$ws.Range("M2").Value = -845.5833
$ws.Range("N2").Value = -1092.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17319860
$ws.Range("I61").Value = 21877562
$ws.Range("J61").Value = 1114695.1
$ws.Range("K61").Value = 21877562
$ws.Range("L61").Value = 1114695.1
$ws.Range("M61").Value = -21877350
$ws.Range("N61").Value = -1115119.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2652.5386
$ws.Range("I63").Value = 3031.5557
$ws.Range("J63").Value = 1799.75
$ws.Range("K63").Value = 3031.5557
$ws.Range("L63").Value = 1799.75
$ws.Range("M63").Value = -2345.5557
$ws.Range("N63").Value = -3171.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2652.5386
$ws.Range("I66").Value = 3031.5557
$ws.Range("J66").Value = 1799.75
$ws.Range("K66").Value = 15157.7785
$ws.Range("L66").Value = 8998.75
$ws.Range("M66").Value = -11725.7785
$ws.Range("N66").Value = -15862.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 916.63635
$ws.Range("I116").Value = 958.5833
$ws.Range("J116").Value = 866.3
$ws.Range("K116").Value = 958.5833
$ws.Range("L116").Value = 866.3
$ws.Range("M116").Value = 1335.4167
$ws.Range("N116").Value = -5454.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2130153
$ws.Range("I132").Value = 1304.5853
$ws.Range("J132").Value = 16677283
$ws.Range("K132").Value = 3913.7559
$ws.Range("L132").Value = 50031849
$ws.Range("M132").Value = -1383.7559
$ws.Range("N132").Value = -50036909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 17319860
$ws.Range("I136").Value = 21877562
$ws.Range("J136").Value = 1114695.1
$ws.Range("K136").Value = 65632686
$ws.Range("L136").Value = 3344085.3
$ws.Range("M136").Value = -65630136
$ws.Range("N136").Value = -3349185.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 916.63635
$ws.Range("I3").Value = 958.5833
$ws.Range("J3").Value = 866.3
$ws.Range("K3").Value = 958.5833
$ws.Range("L3").Value = 866.3
$ws.Range("M3").Value = -844.5833
$ws.Range("N3").Value = -1094.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 621062.1
$ws.Range("I105").Value = 849074.1
$ws.Range("J105").Value = 5429.6
$ws.Range("K105").Value = 849074.1
$ws.Range("L105").Value = 5429.6
$ws.Range("M105").Value = -847327.1
$ws.Range("N105").Value = -8923.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 25000370
$ws.Range("I16").Value = 25000370
$ws.Range("K16").Value = 25000370
$ws.Range("M16").Value = -25000083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13701503
$ws.Range("I31").Value = 20835864
$ws.Range("K31").Value = 20835864
$ws.Range("M31").Value = -20835569

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 13701503
$ws.Range("I34").Value = 20835864
$ws.Range("K34").Value = 20835864
$ws.Range("M34").Value = -20835662

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15994.45
$ws.Range("I99").Value = 8841.143
$ws.Range("J99").Value = 32685.5
$ws.Range("K99").Value = 8841.143
$ws.Range("L99").Value = 32685.5
$ws.Range("M99").Value = -7343.143
$ws.Range("N99").Value = -35681.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 889.7742
$ws.Range("I107").Value = 575.0357
$ws.Range("K107").Value = 575.0357
$ws.Range("M107").Value = 1344.9643

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 25000370
$ws.Range("I113").Value = 25000370
$ws.Range("K113").Value = 25000370
$ws.Range("M113").Value = -24998200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I122").Value = 7250
$ws.Range("J122").Value = 2281.25
$ws.Range("K122").Value = 21750
$ws.Range("L122").Value = 6843.75
$ws.Range("M122").Value = -19300
$ws.Range("N122").Value = -11743.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 15994.45
$ws.Range("I126").Value = 8841.143
$ws.Range("J126").Value = 32685.5
$ws.Range("K126").Value = 26523.429
$ws.Range("L126").Value = 98056.5
$ws.Range("M126").Value = -24053.429
$ws.Range("N126").Value = -102996.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 79.14286
$ws.Range("I2").Value = 75.90000000000001
$ws.Range("J2").Value = 87.25
$ws.Range("K2").Value = 455.4
$ws.Range("L2").Value = 523.5
$ws.Range("M2").Value = -342.4
$ws.Range("N2").Value = -749.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 96999.8
$ws.Range("J37").Value = 96999.8
$ws.Range("L37").Value = 290999.4
$ws.Range("N37").Value = -291223.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4273.7
$ws.Range("I129").Value = 2441.9167
$ws.Range("J129").Value = 7021.375
$ws.Range("K129").Value = 7325.750100000001
$ws.Range("L129").Value = 21064.125
$ws.Range("M129").Value = -2325.750100000001
$ws.Range("N129").Value = -31064.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4232.52
$ws.Range("J131").Value = 5065.7334
$ws.Range("L131").Value = 15197.2002
$ws.Range("N131").Value = -25277.2002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 11542050
$ws.Range("I140").Value = 16667555
$ws.Range("K140").Value = 50002665
$ws.Range("M140").Value = -49997485

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11274919
$ws.Range("I132").Value = 2904.2273
$ws.Range("K132").Value = 8712.6819
$ws.Range("M132").Value = -6182.6819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2637.55
$ws.Range("I132").Value = 1651.5
$ws.Range("K132").Value = 4954.5
$ws.Range("M132").Value = -2424.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3235.818
$ws.Range("I136").Value = 1732.8334
$ws.Range("J136").Value = 5039.4
$ws.Range("K136").Value = 5198.5002
$ws.Range("L136").Value = 15118.2
$ws.Range("M136").Value = -2648.5002
$ws.Range("N136").Value = -20218.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2594.3572
$ws.Range("I81").Value = 2733.3845
$ws.Range("K81").Value = 5466.769
$ws.Range("M81").Value = -4405.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2594.3572
$ws.Range("I84").Value = 2733.3845
$ws.Range("K84").Value = 27333.845
$ws.Range("M84").Value = -22029.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 668977.4
$ws.Range("I136").Value = 1796.9
$ws.Range("K136").Value = 5390.700000000001
$ws.Range("M136").Value = -2840.700000000001
